# Updates the "cryptos" worksheet with refreshed price/volume figures
# (and a row reorder for TheSandbox / InternetComputer(DFINITY)) as
# produced by the scheduled GitHub Actions scraper run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '27.440.72'; E = '  -2.47%  ' },
    @{ Row = 3; D = '1.739.53'; E = '  -3.41%  ' },
    @{ Row = 4; D = '1.004'; E = '  -0.18%  ' },
    @{ Row = 5; D = '324.32'; E = '  -4.11%  ' },
    @{ Row = 6; D = '1.001'; E = '  -0.04%  ' },
    @{ Row = 7; D = '0.4242'; E = '  -9.20%  ' },
    @{ Row = 8; D = '0.3620'; E = '  -2.55%  ' },
    @{ Row = 9; D = '45.07'; E = '  -1.09%  ' },
    @{ Row = 10; D = '1.123'; E = '  -2.10%  ' },
    @{ Row = 11; D = '0.07454'; E = '  -3.58%  ' },
    @{ Row = 12; D = '1.003'; E = '  +0.07%  ' },
    @{ Row = 13; D = '21.59'; E = '  -4.37%  ' },
    @{ Row = 14; D = '6.080'; E = '  -5.37%  ' },
    @{ Row = 15; D = '7.171'; E = '  -2.94%  ' },
    @{ Row = 16; D = '1.736.74'; E = '  -3.51%  ' },
    @{ Row = 17; D = '0.00001067'; E = '  -2.98%  ' },
    @{ Row = 18; D = '86.93'; E = '  +5.65%  ' },
    @{ Row = 19; D = '0.06049'; E = '  -10.47%  ' },
    @{ Row = 20; D = '1.001'; E = '  -0.08%  ' },
    @{ Row = 21; D = '16.83'; E = '  -4.26%  ' },
    @{ Row = 22; D = '6.065' },
    @{ Row = 23; D = '0.5229'; E = '  -5.01%  ' },
    @{ Row = 24; D = '27.446.21'; E = '  -2.47%  ' },
    @{ Row = 25; D = '11.32'; E = '  -5.64%  ' },
    @{ Row = 26; D = '2.409'; E = '  +0.32%  ' },
    @{ Row = 27; D = '20.18'; E = '  -3.57%  ' },
    @{ Row = 28; D = '2.372'; E = '  -1.69%  ' },
    @{ Row = 29; D = '149.44'; E = '  -1.30%  ' },
    @{ Row = 30; D = '1.936.07' },
    @{ Row = 31; D = '1.275'; E = '  +0.45%  ' },
    @{ Row = 32; D = '126.62'; E = '  -6.06%  ' },
    @{ Row = 33; D = '3.734'; E = '  -7.97%  ' },
    @{ Row = 34; D = '5.611'; E = '  -5.97%  ' },
    @{ Row = 35; D = '0.09050'; E = '  -5.88%  ' },
    @{ Row = 36; D = '12.52'; E = '  +2.43%  ' },
    @{ Row = 37; D = '0.2155'; E = '  -3.79%  ' },
    @{ Row = 38; D = '0.06165'; E = '  -3.39%  ' },
    @{ Row = 39; D = '0.02281'; E = '  -4.84%  ' },
    @{ Row = 40; B = 'InternetComputer(DFINITY)'; C = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D = '5.035'; E = '  -5.01%  ' },
    @{ Row = 41; B = 'TheSandbox'; C = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D = '0.6409'; E = '  -4.69%  ' },
    @{ Row = 42; D = '1.185'; E = '  -4.71%  ' },
    @{ Row = 43; D = '1.415'; E = '  -4.91%  ' },
    @{ Row = 44; D = '1.000'; E = '  -0.18%  ' },
    @{ Row = 45; D = '7.831'; E = '  -3.53%  ' },
    @{ Row = 46; D = '13.53'; E = '  -4.90%  ' },
    @{ Row = 47; D = '3.747' },
    @{ Row = 48; D = '0.5856'; E = '  -5.10%  ' },
    @{ Row = 49; D = '125.48'; E = '  -3.56%  ' },
    @{ Row = 50; D = '1.940'; E = '  -6.22%  ' },
    @{ Row = 51; D = '0.06826'; E = '  -4.20%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($u.ContainsKey('B')) {
        $ws.Cells.Item($r, 2).Value = $u.B
    }

    if ($u.ContainsKey('C')) {
        $ws.Cells.Item($r, 3).Value = $u.C
    }

    if ($u.ContainsKey('D')) {
        # Column D holds price strings that often look numeric
        # (e.g. "1.004", "1.000"). Force text formatting before
        # assigning so Excel doesn't silently coerce them to numbers.
        $cellD = $ws.Cells.Item($r, 4)
        $cellD.NumberFormat = '@'
        $cellD.Value = $u.D
    }

    if ($u.ContainsKey('E')) {
        $cellE = $ws.Cells.Item($r, 5)
        $cellE.NumberFormat = '@'
        $cellE.Value = $u.E
    }
}
